$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FPCbS")

# Split "natural gas nonpeaker" into two rows:
#   "natural gas steam turbine" and "natural gas combined cycle"
# Insert a new row below the current "natural gas nonpeaker" row (row 3)
# so the rest of the table shifts down by one.
$ws.Rows.Item(4).Insert()

# Relabel the original row with the first new source name.
$ws.Range("A3").Value = "natural gas steam turbine"
$ws.Range("B3").Value = 0

# Fill the newly inserted row with the second new source name.
$ws.Range("A4").Value = "natural gas combined cycle"
$ws.Range("B4").Value = 0
